$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; this shifts Item / Purchase Date /
# Qty / Unit / cost / cost-per-unit one column to the right (B..G) and
# shifts their formulas (e.g. =E2/C2 -> =F2/D2) automatically.
$ws.Range("A1").EntireColumn.Insert()

# New "Code" column with the blue header fill matching the other headers.
$ws.Range("A1").Value = "Code"
$ws.Range("A1").Interior.Color = 15773696

$ws.Range("A2").Value = "BM00001"
$ws.Range("A3").Value = "BM00002"
$ws.Range("A4").Value = "BM00003"
$ws.Range("A5").Value = "BM00004"

# Match the selection left behind in the saved file.
$ws.Range("A2").Select()
